$wb = $excel.ActiveWorkbook

$sheetInstructions = $wb.Worksheets.Item("instructions")
$sheetData = $wb.Worksheets.Item("Data")
$sheetCat = $wb.Worksheets.Item("Cat")

$sheetInstructions.Range("K2").Value = "Improved Water Source"
$sheetData.Range("L2").Value = "Improved Water Source"
$sheetCat.Range("L2").Value = "Improved Water Source"

$sheetCat.Activate()
